# Fix: "Invoice id is not generated" + a couple of data corrections on the
# sales-invoice template.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Invoice1a")

# --- INVOICE # (F3): was stuck showing a stale/hand-typed invoice number;
# should show the (now properly generated) invoice id "1".
# Written through TEXT()+PasteSpecial(values) so the numeric-looking
# string lands as real text (matching the cell's existing inlineStr/text
# storage) instead of Excel auto-coercing it to a number, while keeping
# the cell's original formatting/style untouched.
$ws.Range("F3").Formula = '=TEXT(1,"0")'
$ws.Range("F3").Copy()
$ws.Range("F3").PasteSpecial(-4163)  # xlPasteValues

# --- Company address (A11): trim the city/state suffix.
$ws.Range("A11").Value = "Bangalore East "

# --- Contact numbers (A13, D16): correct the mistyped phone number.
$ws.Range("A13").Formula = '=TEXT(990019361,"0")'
$ws.Range("A13").Copy()
$ws.Range("A13").PasteSpecial(-4163)

$ws.Range("D16").Formula = '=TEXT(990019361,"0")'
$ws.Range("D16").Copy()
$ws.Range("D16").PasteSpecial(-4163)

$excel.CutCopyMode = 0
